$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Overwrite row 1 (the old header row) in place with numeric placeholders
#    0..11 -- keeps its existing bold/border style (s="1") since we only
#    change the cell values, not the formatting.
for ($c = 1; $c -le 12; $c++) {
    $ws.Cells.Item(1, $c).Value = $c - 1
}

# 2) Insert a brand-new row at position 2 (pushes the old rows 2..66 down
#    to 3..67) to hold the header text that used to live in row 1.
$ws.Rows.Item(2).Insert()

# Excel inherits formatting from the row above on insert; strip that back
# out so the new header-text row carries no style (matches the target).
$ws.Range("A2:L2").ClearFormats()

# 3) Populate the new row 2 with the former header labels.
$ws.Range("A2").Value = "Lg."
$ws.Range("B2").Value = "Threading"
$ws.Range("C2").Value = "Min.Thread Lg."
$ws.Range("D2").Value = "HeadDia."
$ws.Range("E2").Value = "Head Ht."
$ws.Range("F2").Value = "TensileStrength, psi"
$ws.Range("G2").Value = "Specifications Met"
$ws.Range("H2").Value = "Pkg.Qty."
$ws.Range("J2").Value = "Pkg."

# I2, K2, L2 stay empty (I1 was already blank; K1/L1 text is dropped here).
